$wb = $excel.ActiveWorkbook

# --- CT-001 - Adicionar Usuario (C7: Dados de Entrada) ---
$ws = $wb.Worksheets.Item("CT-001")
$h = $ws.Rows("7").RowHeight
$ws.Range("C7").Value = "Nome de usuário para login: Admin`nSenha para login: admin123`nUser Role: Admin`nStatus: Enabled`nEmployee Name: Qualquer um que tenha a letra A`nUsername: Nome Usuario + numero aleatorio(1 a 999)`nPassword:Abcd@1234`nConfirm Password: Abcd@1234"
$ws.Rows("7").RowHeight = $h

# --- CT-002 - Deletar usuario (C7: Dados de Entrada) ---
$ws = $wb.Worksheets.Item("CT-002")
$h = $ws.Rows("7").RowHeight
$ws.Range("C7").Value = "Nome de usuário para login: Admin`nSenha para login: admin123`nUsername: Nome do usuario cadastrado previamente"
$ws.Rows("7").RowHeight = $h

# --- CT-003 - Editar usuario (C7: Dados de Entrada) ---
$ws = $wb.Worksheets.Item("CT-003")
$h = $ws.Rows("7").RowHeight
$ws.Range("C7").Value = "Nome de usuário para login: Admin`nSenha para login: admin123`n[Antigo]Username: Nome do usuário cadastrado previamente`n[Novo]Username: Nome do novo usuario gerado"
$ws.Rows("7").RowHeight = $h

# --- CT-004 - Ordenar usuarios (C4: Procedimentos, C5: Resultado Esperados) ---
$ws = $wb.Worksheets.Item("CT-004")
$h4 = $ws.Rows("4").RowHeight
$h5 = $ws.Rows("5").RowHeight
$ws.Range("C4").Value = "1 - Acessar menu [Admin]`n2 - Clicar no ícone seta da coluna Employee Name`n3 - Selecionar a opção Ascending pelo nome do usuario`n4 - Verificar a ordenação dos registros dos usuários"
$ws.Range("C5").Value = "1 - Tela de gerenciamento de usuários/admin deve ser exibida`n2 - Lista de opções de ordenação deve ser exibida`n3 - `n4 - Registro de usuários deve ser ordenado na ordem alfábetica do Nome do usuario"
$ws.Rows("4").RowHeight = $h4
$ws.Rows("5").RowHeight = $h5

# --- CT-005 - Pesquisar usuario (C7: Dados de Entrada) ---
$ws = $wb.Worksheets.Item("CT-005")
$h = $ws.Rows("7").RowHeight
$ws.Range("C7").Value = "Nome de usuário para login: Admin`nSenha para login: admin123`nUsername: Nome do usuario cadastrado previamente`nUser Role: ESS`nStatus: Admin"
$ws.Rows("7").RowHeight = $h

# --- CT-006 - Resetar pesquisa de usuario (C7: Dados de Entrada) ---
$ws = $wb.Worksheets.Item("CT-006")
$h = $ws.Rows("7").RowHeight
$ws.Range("C7").Value = "Nome de usuário para login: Admin`nSenha para login: admin123`nUsername: Admin`nUser Role: ESS`nStatus: Disabled"
$ws.Rows("7").RowHeight = $h
